$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @(1.02, 1.016575399075721, 1.022036165938839, 1.018079921344615, 1.01492247086263, 1.026428262455471, 1.021794426393935, 1.024871731696536, 1.020927201768244, 1.017779172051364)
  3 = @(1.02, 1.01758499701192, 1.022738956851938, 1.01893769224147, 1.016575723826187, 1.02655242823776, 1.022439219096827, 1.025381386959106, 1.021590562361069, 1.019235123611414)
  4 = @(1.019999999999999, 1.018238117224918, 1.023193354322733, 1.019492974428462, 1.017645160628615, 1.026631240251872, 1.022855740859572, 1.025710143057602, 1.022019413185946, 1.020176446575417)
  5 = @(1.02, 1.018512652216644, 1.023384297700813, 1.01972647432373, 1.018094680237479, 1.026664006033661, 1.023030678983413, 1.025848106794078, 1.022199609524285, 1.020571999226914)
  6 = @(1.02, 1.018558745720429, 1.023416352922965, 1.019765683462125, 1.018170152560844, 1.026669486040007, 1.023060042048678, 1.025871257127548, 1.022229859863914, 1.020638403905121)
  7 = @(1.02, 1.018241785718215, 1.023195906053109, 1.019496094232499, 1.017651167404378, 1.026631679511256, 1.022858079049406, 1.025711987499738, 1.022021821341536, 1.020181732667243)
  8 = @(1.02, 1.016916630497611, 1.02227375081335, 1.018369757953472, 1.015481267661749, 1.026470541901546, 1.022012482806205, 1.025044184327718, 1.021151468073679, 1.018271381809822)
  9 = @(1.02, 1.014580293964472, 1.020646085865126, 1.016386895540633, 1.011654798108463, 1.026174872547915, 1.02051703326557, 1.023859577264158, 1.019614809869173, 1.014898893402669)
  10 = @(1.02, 1.013021841066643, 1.019559166275683, 1.015066243486943, 1.009101476840297, 1.025969883242132, 1.019516401444701, 1.023064553125041, 1.018588341998398, 1.012646016316115)
  11 = @(1.02, 1.012346787360823, 1.019088090513704, 1.014494680898912, 1.007995207076624, 1.02587925339289, 1.01908223903227, 1.022719042840219, 1.018143383010055, 1.011669331470965)
  12 = @(1.02, 1.012096005985239, 1.018913047064612, 1.014282420028335, 1.007584180848626, 1.025845308785465, 1.01892083842139, 1.022590515406041, 1.017978030978306, 1.011306363467904)
  13 = @(1.02, 1.012149801108823, 1.018950597416832, 1.014327948770025, 1.007672352411251, 1.025852602718548, 1.018955465448314, 1.022618093558412, 1.018013502929098, 1.01138422977092)
  14 = @(1.02, 1.012326058443739, 1.019073622692163, 1.014477134475953, 1.007961233780678, 1.02587645324605, 1.019068900332674, 1.022708422590504, 1.018129716481771, 1.011639332214027)
  15 = @(1.02, 1.012434651545663, 1.019149414020556, 1.014569058388275, 1.008139208510017, 1.02589111115538, 1.019138773612734, 1.022764052164874, 1.018201309603701, 1.011796484677385)
  16 = @(1.02, 1.013066637129772, 1.019590420914886, 1.015104182255978, 1.009174881589792, 1.025975858689142, 1.019545196717251, 1.023087456941997, 1.018617862052851, 1.012710810306764)
  17 = @(1.02, 1.013463001596714, 1.019866937144499, 1.01543992826175, 1.009824347977968, 1.026028518441906, 1.019799898423503, 1.023289982784614, 1.018879022363394, 1.013284023073536)
  18 = @(1.02, 1.013694171861702, 1.020028182660723, 1.015635791087851, 1.010203107162639, 1.026059053736506, 1.019948376552272, 1.023407991170436, 1.01903130535936, 1.013618255584452)
  19 = @(1.02, 1.013772991184794, 1.020083156136956, 1.015702579988222, 1.010332243731739, 1.026069434909006, 1.019998989352903, 1.023448208378639, 1.019083221904484, 1.013732201324003)
  20 = @(1.02, 1.013420477769123, 1.019837273885157, 1.015403903013792, 1.009754673016978, 1.02602288719253, 1.019772580145495, 1.023268266248874, 1.018851007240016, 1.013222534488014)
  21 = @(1.02, 1.012274156066409, 1.019037396607302, 1.014433201822409, 1.007876168477913, 1.025869437607231, 1.019035500274498, 1.022681828176307, 1.018095496547077, 1.011564216034856)
  22 = @(1.02, 1.011553207013932, 1.018534106318378, 1.013823131584085, 1.006694446813995, 1.025771334037443, 1.018571296915223, 1.022312014567588, 1.017620045934155, 1.01052049852883)
  23 = @(1.02, 1.011935416200788, 1.018800945628863, 1.014146517982557, 1.007320962519004, 1.025823494501011, 1.01881745334963, 1.02250816381488, 1.01787213235877, 1.011073896521112)
  24 = @(1.02, 1.013439692525261, 1.019850677563984, 1.015420181186343, 1.009786156327183, 1.026025432268605, 1.019784924361976, 1.023278079391863, 1.018863666214571, 1.013250318879863)
  25 = @(1.02, 1.015184446388862, 1.021067196074201, 1.016899290678068, 1.012644412586164, 1.026252698885156, 1.020904286338894, 1.024166757786121, 1.020012428806697, 1.015771536848917)
}

$cols = @("B","C","D","E","F","I","J","K","L","M")

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])$row").Value = $vals[$i]
  }
}